# day 1 R part 21
#
# Changes applied (per the target diff):
#  1. In the code textbox ("TextBox 3") inside Group 20 / Group 14, the
#     word "Iris" that began the first line ("Iris %>%") is replaced with
#     "MyFavoriteIrisFlowers", turning the line into
#     "MyFavoriteIrisFlowers %>%" (as two runs, matching how PowerPoint
#     splits a run when you type new text immediately before existing text).
#  2. The "Aesthetic map" textbox ("TextBox 7") is moved to the right.
#  3. The arrow connector that starts at "Aesthetic map" ("Straight Arrow
#     Connector 8") is moved along with it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($container, $name) {
    for ($i = 1; $i -le $container.Count; $i++) {
        $sh = $container.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# The top-level shape on the slide is the big "Group 20" group; all the
# shapes we need to touch live inside it (PowerPoint flattens the nested
# "Group 14" wrapper into Group 20's GroupItems).
$g20 = $s.Shapes.Item(1)

# --- 1. Text edit: "Iris %>%" -> "MyFavoriteIrisFlowers %>%" -----------
$tb3 = Get-ShapeByName $g20.GroupItems "TextBox 3"
$tr = $tb3.TextFrame.TextRange
# "Iris" is exactly the first 4 characters of the textbox's text.
$irisRange = $tr.Characters(1, 4)
$irisRange.Text = "MyFavoriteIrisFlowers"

# --- 2 & 3. Reposition "Aesthetic map" textbox and its arrow -----------
# Target local offsets (EMU), taken directly from the XML diff:
#   TextBox 7: off x=6223380, y=1450655  (was x=2757999, y=1492245)
#   Connector8: off x=6223390, y=1973875 (was x=2758009, y=2015465)
# The host's Shape.Left/Top setters write points*12700 straight into the
# shape's local <a:off> (no re-projection through the parent group's
# child-offset transform), so feeding target_EMU/12700 reproduces the
# exact local offsets from the diff.
$tb7 = Get-ShapeByName $g20.GroupItems "TextBox 7"
$tb7.Left = 490.0299225598428
$tb7.Top = 114.22480394960625

$conn8 = Get-ShapeByName $g20.GroupItems "Straight Arrow Connector 8"
$conn8.Left = 490.030715961419
$conn8.Top = 155.42322834645668
